# Add a new localization row ("findProduct") to the "login" sheet.
# Mirrors the source xlsx diff: a new row 9 with Key/en/id/es columns,
# and the sheet selection left on D10 (the cell right below/after the
# freshly entered data), as Excel would leave it after typing the last
# value and pressing Enter/Tab off the table.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("login")

$ws1.Range("A9").Value = "findProduct"
$ws1.Range("B9").Value = "Find Product"
$ws1.Range("C9").Value = "Cari Product"
$ws1.Range("D9").Value = "selemente"

# Leave the selection where Excel would land after entering the row.
$ws1.Range("D10").Select()
